{"js": "// Fix a typo (\"Unusable\" -> \"Unuseable\") and make dash usage consistent\n// by converting every en dash (U+2013) used as an em dash (U+2014)\n// throughout the document body, per the commit message\n// \"consistency of grammar, shortenings\".\n\nconst body = context.document.body;\n\nconst enDash = String.fromCharCode(0x2013);\nconst emDash = String.fromCharCode(0x2014);\n\n// 1) Typo fix: \"Unusable\" -> \"Unuseable\" (appears once, in the Theme 2 /\n// Table 5.7 caption line).\nconst typoResults = body.search(\"Unusable\", { matchCase: true, matchWholeWord: false });\ntypoResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < typoResults.items.length; i++) {\n  typoResults.items[i].insertText(\"Unuseable\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Dash consistency: replace every en dash with an em dash (document-wide).\nconst dashResults = body.search(enDash, { matchCase: true });\ndashResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < dashResults.items.length; i++) {\n  dashResults.items[i].insertText(emDash, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fix a typo (\"Unusable\" -> \"Unuseable\") and make dash usage consistent by\n# converting every en dash (U+2013) used as an em dash (U+2014) throughout\n# the document, per the commit message \"consistency of grammar, shortenings\".\n\n$d = $word.ActiveDocument\n\n$enDash = [char]0x2013\n$emDash = [char]0x2014\n\n# 1) Typo fix: \"Unusable\" -> \"Unuseable\" (Table 5.7 / Theme 2 caption).\n$find = $d.Content.Find\n$find.Execute(\"Unusable\", $true, $false, $false, $false, $false, $true, 1, $false, \"Unuseable\", 2)\n\n# 2) Dash consistency: replace every en dash with an em dash (document-wide).\n$find2 = $d.Content.Find\n$find2.Execute($enDash, $true, $false, $false, $false, $false, $true, 1, $false, $emDash, 2)\n"}
